$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking strings
# (e.g. "219.32", "0.06416") keep their exact literal representation
# instead of being parsed into floating-point numbers.
$ws.Range("D2:D50").NumberFormat = "@"

$ws.Range('D2').Value = '26.344.67'
$ws.Range('E2').Value = '  +1.09%  '
$ws.Range('D3').Value = '1.666.81'
$ws.Range('E3').Value = '  +0.93%  '
$ws.Range('E4').Value = '  +0.91%  '
$ws.Range('D5').Value = '219.32'
$ws.Range('E5').Value = '  +0.79%  '
$ws.Range('D6').Value = '0.5352'
$ws.Range('E6').Value = '  +1.89%  '
$ws.Range('E7').Value = '  +0.85%  '
$ws.Range('D8').Value = '0.2661'
$ws.Range('E8').Value = '  +2.52%  '
$ws.Range('D9').Value = '0.06416'
$ws.Range('E9').Value = '  +1.57%  '
$ws.Range('D10').Value = '20.77'
$ws.Range('E10').Value = '  +2.07%  '
$ws.Range('D11').Value = '0.07844'
$ws.Range('E11').Value = '  +0.78%  '
$ws.Range('D12').Value = '4.568'
$ws.Range('E12').Value = '  +1.48%  '
$ws.Range('D13').Value = '1.664.39'
$ws.Range('E13').Value = '  +0.95%  '
$ws.Range('D14').Value = '1.893.05'
$ws.Range('E14').Value = '  +0.78%  '
$ws.Range('D15').Value = '0.5542'
$ws.Range('E15').Value = '  +1.12%  '
$ws.Range('D16').Value = '0.0₅8201'
$ws.Range('E16').Value = '  +0.03%  '
$ws.Range('D17').Value = '65.83'
$ws.Range('E17').Value = '  +0.50%  '
$ws.Range('D18').Value = '26.359.11'
$ws.Range('E18').Value = '  +1.12%  '
$ws.Range('E19').Value = '  +0.93%  '
$ws.Range('D20').Value = '4.691'
$ws.Range('E20').Value = '  +2.51%  '
$ws.Range('D21').Value = '193.47'
$ws.Range('E21').Value = '  +1.28%  '
$ws.Range('D22').Value = '10.30'
$ws.Range('E22').Value = '  +2.16%  '
$ws.Range('D23').Value = '6.050'
$ws.Range('E23').Value = '  +0.34%  '
$ws.Range('E24').Value = '  +0.87%  '
$ws.Range('D25').Value = '146.46'
$ws.Range('E25').Value = '  +2.29%  '
$ws.Range('D26').Value = '0.1230'
$ws.Range('E26').Value = '  -0.46%  '
$ws.Range('D27').Value = '7.220'
$ws.Range('E27').Value = '  +0.03%  '
$ws.Range('D28').Value = '16.12'
$ws.Range('E28').Value = '  +0.58%  '
$ws.Range('D29').Value = '1.497'
$ws.Range('E29').Value = '  +4.92%  '
$ws.Range('D30').Value = '0.05884'
$ws.Range('E30').Value = '  +1.32%  '
$ws.Range('D31').Value = '1.289'
$ws.Range('E31').Value = '  +1.19%  '
$ws.Range('D32').Value = '3.636'
$ws.Range('E32').Value = '  +2.42%  '
$ws.Range('D33').Value = '3.288'
$ws.Range('E33').Value = '  +0.73%  '
$ws.Range('D34').Value = '1.609'
$ws.Range('E34').Value = '  +1.73%  '
$ws.Range('D35').Value = '0.9714'
$ws.Range('E35').Value = '  +2.80%  '
$ws.Range('D36').Value = '2.829'
$ws.Range('E36').Value = '  +1.76%  '
$ws.Range('D37').Value = '2.420'
$ws.Range('E37').Value = '  +0.39%  '
$ws.Range('D38').Value = '0.5839'
$ws.Range('E38').Value = '  +1.77%  '
$ws.Range('D39').Value = '0.01605'
$ws.Range('E39').Value = '  -0.22%  '
$ws.Range('D40').Value = '0.8706'
$ws.Range('E40').Value = '  +3.40%  '
$ws.Range('D41').Value = '5.851'
$ws.Range('E41').Value = '  +1.88%  '
$ws.Range('D42').Value = '1.054.56'
$ws.Range('E42').Value = '  +2.43%  '
$ws.Range('D43').Value = '105.19'
$ws.Range('E43').Value = '  +1.59%  '
$ws.Range('E44').Value = '  +0.88%  '
$ws.Range('D45').Value = '1.804.46'
$ws.Range('E45').Value = '  +0.58%  '
$ws.Range('D46').Value = '57.86'
$ws.Range('E46').Value = '  +1.70%  '
$ws.Range('E47').Value = '  -4.94%  '
$ws.Range('D48').Value = '1.013'
$ws.Range('E48').Value = '  +1.02%  '
$ws.Range('D49').Value = '0.4388'
$ws.Range('E49').Value = '  +1.47%  '
$ws.Range('D50').Value = '7.990'
$ws.Range('E50').Value = '  +1.64%  '

# Restore the default style on column D (remove the temporary text format)
$ws.Range("D2:D50").Style = "Normal"
